# Update countries & provincias Spain
# Applies the 10-Abril-2020 17:52 -> 18:22 data refresh to the "paises" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 18:22"

# --- Country names that moved rank (ranking reshuffled by updated totals) ---
# Reino Unido overtakes Iran -> they swap row positions 10/11
$ws.Range("A10").Value = "Reino Unido"
$ws.Range("A11").Value = "Iran"

# Argelia overtakes Egipto & Islandia (ranks 61/62/63 shuffle)
$ws.Range("A57").Value = "Argelia"
$ws.Range("A58").Value = "Egipto"
$ws.Range("A59").Value = "Islandia"

# --- Row 4 (Estados Unidos) updated counts -----------------------------
$ws.Range("B4").Value = 478004
$ws.Range("C4").Value = 9438
$ws.Range("D4").Value = 26098
$ws.Range("E4").Value = 433995
$ws.Range("F4").Value = 10571
$ws.Range("G4").Value = 1220
$ws.Range("H4").Value = 17911

# --- Row 6 (Italia) updated counts -------------------------------------
$ws.Range("B6").Value = 147577
$ws.Range("C6").Value = 3951
$ws.Range("D6").Value = 30455
$ws.Range("E6").Value = 98273
$ws.Range("F6").Value = 3497
$ws.Range("G6").Value = 570
$ws.Range("H6").Value = 18849

# --- Row 7 (Alemania) updated counts ------------------------------------
$ws.Range("B7").Value = 119624
$ws.Range("C7").Value = 1389
$ws.Range("E7").Value = 64610

# --- Row 10 (now Reino Unido) updated counts ----------------------------
$ws.Range("B10").Value = 70272
$ws.Range("C10").Value = 5195
$ws.Range("D10").Value = 135
$ws.Range("E10").Value = 61206
$ws.Range("F10").Value = 1559
$ws.Range("G10").Value = 953
$ws.Range("H10").Value = 8931

# --- Row 11 (now Iran) updated counts -----------------------------------
$ws.Range("B11").Value = 68192
$ws.Range("C11").Value = 1972
$ws.Range("D11").Value = 35465
$ws.Range("E11").Value = 28495
$ws.Range("F11").Value = 3969
$ws.Range("G11").Value = 122
$ws.Range("H11").Value = 4232

# --- Row 14 (Suiza) updated counts ------------------------------------
$ws.Range("B14").Value = 24548
$ws.Range("C14").Value = 497
$ws.Range("E14").Value = 12947
$ws.Range("G14").Value = 53
$ws.Range("H14").Value = 1001

# --- Row 16 (Paises Bajos) updated counts -------------------------------
$ws.Range("D16").Value = 5580
$ws.Range("E16").Value = 15132
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = 531

# --- Row 17 (Canada) updated counts -------------------------------------
$ws.Range("B17").Value = 18397
$ws.Range("C17").Value = 252
$ws.Range("E17").Value = 17250
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 974

# --- Row 19 (Austria) updated counts ------------------------------------
$ws.Range("B19").Value = 13531
$ws.Range("C19").Value = 287
$ws.Range("E19").Value = 7148

# --- Row 31 updated counts -----------------------------------------------
$ws.Range("B31").Value = 5674
$ws.Range("C31").Value = 105
$ws.Range("D31").Value = 346
$ws.Range("E31").Value = 5209
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = 119

# --- Row 42 updated counts -------------------------------------------------
$ws.Range("B42").Value = 3223
$ws.Range("C42").Value = 108
$ws.Range("E42").Value = 2669
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 54

# --- Row 53 updated count ---------------------------------------------------
$ws.Range("F53").Value = 77

# --- Row 57 (now Argelia) updated counts -------------------------------------
$ws.Range("B57").Value = 1761
$ws.Range("C57").Value = 95
$ws.Range("D57").Value = 405
$ws.Range("E57").Value = 1100
$ws.Range("F57").Value = 46
$ws.Range("G57").Value = 21
$ws.Range("H57").Value = 256

# --- Row 58 (now Egipto) updated counts ---------------------------------------
$ws.Range("B58").Value = 1699
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 348
$ws.Range("E58").Value = 1233
$ws.Range("F58").Value = 0
$ws.Range("H58").Value = 118

# --- Row 59 (now Islandia) updated counts -------------------------------------
$ws.Range("B59").Value = 1675
$ws.Range("C59").Value = 27
$ws.Range("D59").Value = 751
$ws.Range("E59").Value = 918
$ws.Range("F59").Value = 11
$ws.Range("H59").Value = 6

# --- Row 87 updated count -------------------------------------------------------
$ws.Range("F87").Value = 11

# --- Row 100 updated counts ------------------------------------------------------
$ws.Range("D100").Value = 170
$ws.Range("E100").Value = 195
